$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header for United States in column F
$ws.Range("F1").Value = "United States"

# Updated values after 85% data collected
$ws.Range("B2").Value = 1.1461038961039
$ws.Range("C2").Value = 1.13821138211382
$ws.Range("D2").Value = 1.14448669201521
$ws.Range("E2").Value = 0.978461538461538

$ws.Range("B3").Value = 0.331168831168831
$ws.Range("C3").Value = 0.493224932249323
$ws.Range("D3").Value = 0.555133079847909
$ws.Range("E3").Value = 0.498461538461538

$ws.Range("B4").Value = 0.227272727272727
$ws.Range("C4").Value = 0.514905149051491
$ws.Range("D4").Value = 0.414448669201521
$ws.Range("E4").Value = 0.326153846153846

$ws.Range("B5").Value = 0.941558441558442
$ws.Range("C5").Value = 0.905149051490515
$ws.Range("D5").Value = 0.707224334600761
$ws.Range("E5").Value = 0.858461538461538

$ws.Range("B6").Value = 1.13961038961039
$ws.Range("C6").Value = 1.04336043360434
$ws.Range("D6").Value = 1.12167300380228
$ws.Range("E6").Value = 1

$ws.Range("B7").Value = 0.707792207792208
$ws.Range("C7").Value = 0.485094850948509
$ws.Range("D7").Value = 0.760456273764259
$ws.Range("E7").Value = 0.461538461538462

$ws.Range("B8").Value = 0.922077922077922
$ws.Range("C8").Value = 0.921409214092141
$ws.Range("D8").Value = 1.06083650190114
$ws.Range("E8").Value = 0.855384615384615

$ws.Range("B9").Value = 0.931818181818182
$ws.Range("C9").Value = 0.569105691056911
$ws.Range("D9").Value = 0.91254752851711
$ws.Range("E9").Value = 0.6

$ws.Range("B10").Value = 0.844155844155844
$ws.Range("C10").Value = 0.734417344173442
$ws.Range("D10").Value = 0.996197718631179
$ws.Range("E10").Value = 0.935384615384615

$ws.Range("B11").Value = 0.665584415584416
$ws.Range("C11").Value = 0.411924119241192
$ws.Range("D11").Value = 0.722433460076046
$ws.Range("E11").Value = 0.704615384615385

$ws.Range("B12").Value = 1.1525974025974
$ws.Range("C12").Value = 0.520325203252033
$ws.Range("D12").Value = 1.24714828897338
$ws.Range("E12").Value = 1.25846153846154
